# Update cryptocurrency price and volume(1h) figures with the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.007.65"
$ws.Range("E2").Value = "  +0.11%  "

$ws.Range("D3").Value = "1.829.10"
$ws.Range("E3").Value = "  +0.50%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -0.49%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.50"
$ws.Range("E5").Value = "  +0.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  -0.38%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4631"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3708"
$ws.Range("E8").Value = "  +1.81%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07342"
$ws.Range("E9").Value = "  +0.76%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8773"
$ws.Range("E10").Value = "  +1.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07896"
$ws.Range("E11").Value = "  +3.54%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.76"
$ws.Range("E12").Value = "  -0.37%  "

$ws.Range("D13").Value = "1.808.40"
$ws.Range("E13").Value = "  -1.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.341"
$ws.Range("E14").Value = "  +0.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.563"
$ws.Range("E15").Value = "  +1.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.48"
$ws.Range("E16").Value = "  -1.91%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.006"
$ws.Range("E17").Value = "  -0.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008848"
$ws.Range("E18").Value = "  +2.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("E19").Value = "  -0.35%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.79"
$ws.Range("E20").Value = "  +2.13%  "

$ws.Range("D21").Value = "27.025.45"
$ws.Range("E21").Value = "  -1.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.106"
$ws.Range("E22").Value = "  -1.07%  "

$ws.Range("E23").Value = "  -0.52%  "

$ws.Range("D24").Value = "2.071.61"
$ws.Range("E24").Value = "  -1.74%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.30"
$ws.Range("E25").Value = "  +0.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.846"
$ws.Range("E26").Value = "  -0.81%  "

$ws.Range("E27").Value = "  +1.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.042"
$ws.Range("E28").Value = "  -2.90%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.112"
$ws.Range("E29").Value = "  +0.45%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.72"
$ws.Range("E30").Value = "  -0.17%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08893"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.961"
$ws.Range("E32").Value = "  +0.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7306"
$ws.Range("E33").Value = "  -0.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.441"
$ws.Range("E34").Value = "  +0.19%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.133"
$ws.Range("E35").Value = "  -0.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.461"
$ws.Range("E36").Value = "  -2.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.078"
$ws.Range("E37").Value = "  +0.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01947"
$ws.Range("E38").Value = "  +1.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05219"
$ws.Range("E39").Value = "  -0.99%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.955"
$ws.Range("E40").Value = "  +0.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.113"
$ws.Range("E41").Value = "  -0.48%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5167"
$ws.Range("E42").Value = "  -0.99%  "

$ws.Range("E43").Value = "  -0.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.172"
$ws.Range("E44").Value = "  -0.87%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4829"
$ws.Range("E45").Value = "  -0.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.005"
$ws.Range("E46").Value = "  -0.38%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.21"
$ws.Range("E47").Value = "  +1.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.22"
$ws.Range("E48").Value = "  -0.92%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.626"
$ws.Range("E49").Value = "  -0.49%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06205"
$ws.Range("E50").Value = "  -0.30%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.76"
$ws.Range("E51").Value = "  +0.26%  "
